# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / handoff / handback timestamp
# columns with the newly generated report's timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the c9315bc4... file
$wsOverview.Range("G2").Value = "2016-08-15 15:03:11"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the c9315bc4... handoff/handback pair
$wsZhCn.Range("H2").Value = "2016-08-15 15:03:01"
$wsZhCn.Range("K2").Value = "2016-08-15 15:03:30"

# de-de sheet: Correspond Handoff Datetime (shares same generate time as Overview)
# and Correspond Handback DateTime for the c9315bc4... handoff/handback pair
$wsDeDe.Range("H2").Value = "2016-08-15 15:03:11"
$wsDeDe.Range("K2").Value = "2016-08-15 15:03:37"
